$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows of data (ID, NombreCompleto, Fallecido). The "false" text must be
# forced as literal text (leading apostrophe) so it is stored as a shared
# string instead of being auto-coerced into a boolean by Excel.
$row2 = New-Object 'object[,]' 1,3
$row2[0,0] = 1
$row2[0,1] = "Persona 1"
$row2[0,2] = "'false"
$ws.Range("A2:C2").Value = $row2

$row3 = New-Object 'object[,]' 1,3
$row3[0,0] = 2
$row3[0,1] = "Persona 2"
$row3[0,2] = "'false"
$ws.Range("A3:C3").Value = $row3

$row4 = New-Object 'object[,]' 1,3
$row4[0,0] = 3
$row4[0,1] = "Tia Abuela Gladis"
$row4[0,2] = "'false"
$ws.Range("A4:C4").Value = $row4
